# Generate Report for Handoff
# Rename the two original e2e fixtures (one .md file + two .png dependents)
# into the new caller/callee markdown fixture set, and add the 4th row
# (callerMd2.md) that the new dependency graph introduces, on all three
# sheets: Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook

$hyperlinkColor = 15570276   # BGR for FF6495ED, matches the workbook's existing custom "HyperLink" font
$xlUnderlineStyleSingle = 2

function Style-AsLink($rng) {
    $f = $rng.Font
    $f.Color = $hyperlinkColor
    $f.Underline = $xlUnderlineStyleSingle
}

# ---------------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de | Latest Handoff Date
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

if ($ov.Hyperlinks.Count() -gt 0) { $ov.Hyperlinks.Delete() }

$ov.Range("A2").Value = "calleeMd1.md"
$ov.Range("B2").Value = "Ready for handoff"
$ov.Range("C2").Value = "Ready for handoff"
$ov.Range("D2").Value = "2016-03-23 07:22:55"

$ov.Range("A3").Value = "calleeMd2.md"
$ov.Range("B3").Value = "Ready for handoff"
$ov.Range("C3").Value = "Ready for handoff"
$ov.Range("D3").Value = "2016-03-23 07:22:55"

$ov.Range("A4").Value = "callerMd1.md"
$ov.Range("B4").Value = "Ready for handoff"
$ov.Range("C4").Value = "Ready for handoff"
$ov.Range("D4").Value = "2016-03-23 07:22:55"

$ov.Range("A5").Value = "callerMd2.md"
$ov.Range("B5").Value = "Ready for handoff"
$ov.Range("C5").Value = "Ready for handoff"
$ov.Range("D5").Value = "2016-03-23 07:22:55"

$ov.Range("D2:D5").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ov.Hyperlinks.Add($ov.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/9bfed87a88e5de207a32c2d38fdd20fbc3132aa6/e2e/calleeMd1.md", "", "", "calleeMd1.md")
$ov.Hyperlinks.Add($ov.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/9bfed87a88e5de207a32c2d38fdd20fbc3132aa6/e2e/calleeMd2.md", "", "", "calleeMd2.md")
$ov.Hyperlinks.Add($ov.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/9bfed87a88e5de207a32c2d38fdd20fbc3132aa6/e2e/callerMd1.md", "", "", "callerMd1.md")
$ov.Hyperlinks.Add($ov.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/9bfed87a88e5de207a32c2d38fdd20fbc3132aa6/e2e/callerMd2.md", "", "", "callerMd2.md")

Style-AsLink $ov.Range("A2:A5")

# ---------------------------------------------------------------------
# Per-locale sheets ("zh-cn" and "de-de") share the same column layout:
# A Source File Name, B File Extension, C Status, D Latest Handoff File,
# E Latest Handoff Datetime, H Latest Handback DateTime,
# I Reference Tokens, J Handoff Reason, K Dependency From
# ---------------------------------------------------------------------
$locales = @(
    @{ Name = "zh-cn"; Xlf = "zh-cn"; HandoffDate = "2016-03-23 07:22:47"; XlfUrlBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2687f762eddd7dfa0346d62cbacebf70b13e322f/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht" },
    @{ Name = "de-de"; Xlf = "de-de"; HandoffDate = "2016-03-23 07:22:55"; XlfUrlBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ec0f1ce7b7f2a6890288e076d341d78d9ee28ab8/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht" }
)

$mdUrlBase = "https://github.com/OpenLocalizationTest/oltest/blob/9bfed87a88e5de207a32c2d38fdd20fbc3132aa6/e2e"

foreach ($loc in $locales) {
    $ws = $wb.Worksheets.Item($loc.Name)
    if ($ws.Hyperlinks.Count() -gt 0) { $ws.Hyperlinks.Delete() }

    $xlfNames = @{
        "calleeMd1" = "calleeMd1.e8f5ecec2b522eb147a4ff0ca19ca72e17f2186d.$($loc.Xlf).xlf"
        "calleeMd2" = "calleeMd2.63b76063f058ecc63ff1dda71ea2a67db72ae6e1.$($loc.Xlf).xlf"
        "callerMd1" = "callerMd1.a3bf9f4e7fa2750ec06df0b78a76ae5cafa0e0fd.$($loc.Xlf).xlf"
        "callerMd2" = "callerMd2.c7d976edeb9cd5406eae7aba4c05d6d92e81ae95.$($loc.Xlf).xlf"
    }

    # Row 2: calleeMd1.md
    $ws.Range("A2").Value = "calleeMd1.md"
    $ws.Range("B2").Value = ".md"
    $ws.Range("C2").Value = "Ready for handoff"
    $ws.Range("D2").Value = $xlfNames["calleeMd1"]
    $ws.Range("E2").Value = $loc.HandoffDate
    $ws.Range("H2").Value = "0001-01-01 00:00:00"
    $ws.Range("I2").ClearContents()
    $ws.Range("J2").Value = "Include"
    $ws.Range("K2").Value = "e2e\callerMd2.md,`ne2e\callerMd1.md"

    # Row 3: calleeMd2.md
    $ws.Range("A3").Value = "calleeMd2.md"
    $ws.Range("B3").Value = ".md"
    $ws.Range("C3").Value = "Ready for handoff"
    $ws.Range("D3").Value = $xlfNames["calleeMd2"]
    $ws.Range("E3").Value = $loc.HandoffDate
    $ws.Range("H3").Value = "0001-01-01 00:00:00"
    $ws.Range("I3").ClearContents()
    $ws.Range("J3").Value = "Include"
    $ws.Range("K3").Value = "e2e\callerMd1.md"

    # Row 4: callerMd1.md
    $ws.Range("A4").Value = "callerMd1.md"
    $ws.Range("B4").Value = ".md"
    $ws.Range("C4").Value = "Ready for handoff"
    $ws.Range("D4").Value = $xlfNames["callerMd1"]
    $ws.Range("E4").Value = $loc.HandoffDate
    $ws.Range("H4").Value = "0001-01-01 00:00:00"
    $ws.Range("I4").Value = "e2e\calleeMd1.md,`ne2e\calleeMd2.md"
    $ws.Range("J4").Value = "Include"
    $ws.Range("K4").ClearContents()

    # Row 5: callerMd2.md (new row)
    $ws.Range("A5").Value = "callerMd2.md"
    $ws.Range("B5").Value = ".md"
    $ws.Range("C5").Value = "Ready for handoff"
    $ws.Range("D5").Value = $xlfNames["callerMd2"]
    $ws.Range("E5").Value = $loc.HandoffDate
    $ws.Range("H5").Value = "0001-01-01 00:00:00"
    $ws.Range("I5").Value = "e2e\calleeMd1.md"
    $ws.Range("J5").Value = "Include"

    $ws.Range("E2:E5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
    $ws.Range("H2:H5").NumberFormat = "yyyy-mm-dd HH:mm:ss"

    $ws.Hyperlinks.Add($ws.Range("A2"), "$mdUrlBase/calleeMd1.md", "", "", "calleeMd1.md")
    $ws.Hyperlinks.Add($ws.Range("D2"), "$($loc.XlfUrlBase)/$($xlfNames['calleeMd1'])", "", "", $xlfNames["calleeMd1"])
    $ws.Hyperlinks.Add($ws.Range("A3"), "$mdUrlBase/calleeMd2.md", "", "", "calleeMd2.md")
    $ws.Hyperlinks.Add($ws.Range("D3"), "$($loc.XlfUrlBase)/$($xlfNames['calleeMd2'])", "", "", $xlfNames["calleeMd2"])
    $ws.Hyperlinks.Add($ws.Range("A4"), "$mdUrlBase/callerMd1.md", "", "", "callerMd1.md")
    $ws.Hyperlinks.Add($ws.Range("D4"), "$($loc.XlfUrlBase)/$($xlfNames['callerMd1'])", "", "", $xlfNames["callerMd1"])
    $ws.Hyperlinks.Add($ws.Range("A5"), "$mdUrlBase/callerMd2.md", "", "", "callerMd2.md")
    $ws.Hyperlinks.Add($ws.Range("D5"), "$($loc.XlfUrlBase)/$($xlfNames['callerMd2'])", "", "", $xlfNames["callerMd2"])

    Style-AsLink $ws.Range("A2:A5")
    Style-AsLink $ws.Range("D2:D5")
}

Write-Output "edit complete"
